$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1. "Last compiled ..." Date paragraph: re-split the two existing runs
#    ("L" / "ast compiled ...") into three runs with different text
#    breakpoints, add sz/szCs=20 (10pt) throughout, and add gramStart/
#    gramEnd proofErr markers around the middle run.
# -----------------------------------------------------------------------
$dateRange = $d.Paragraphs(2).Range

$bodyFrag = '<w:p><w:pPr><w:pStyle w:val="Date"/><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Last compiled Mon. 2021-05-</w:t></w:r>' + `
  '<w:proofErr w:type="gramStart"/>' + `
  '<w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>31,  7</w:t></w:r>' + `
  '<w:proofErr w:type="gramEnd"/>' + `
  '<w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>:37 PM</w:t></w:r>' + `
  '</w:p>'

$pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
  '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  $bodyFrag + `
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$dateRange.InsertXML($pkg)

# InsertXML cannot carry an rStyle reference through, so re-apply the
# "Verbatim Char" character style across the whole (re-inserted) run text.
$dateRange2 = $d.Paragraphs(2).Range
$dateTextLen = $dateRange2.Text.Length
$dateTextRange = $d.Range($dateRange2.Start, $dateRange2.Start + $dateTextLen - 1)
$dateTextRange.Style = "Verbatim Char"

# -----------------------------------------------------------------------
# 2. Body paragraphs: consolidate split runs back into single runs and
#    switch the paragraph style from "First Paragraph" to "Body Text".
# -----------------------------------------------------------------------

# Paragraph: "An ideal infield surface ..." (3 runs -> 1 run)
$merge3 = "An ideal infield surface allows players" + [char]0x2019 + " cleats to penetrate the soil and provide adequate traction but impart minimal disruption during play. This state has been termed the " + [char]0x201C + "cleat-in/cleat-out effect." + [char]0x201D + " As an infield soil dries, it transitions from a cleat-in/cleat-out state to a more brittle condition in which the primary yield mode is chip-forming or clod-forming failure. Large surface irregularities formed in this state may deflect batted balls and induce fielding errors or injuries."
$d.Content.Find.Execute($merge3, $true, $false, $false, $false, $false, $true, 1, $false, $merge3, 2) | Out-Null

# Paragraph: "The goal of this research ..." (2 runs -> 1 run, before the inline math)
$merge4 = "The goal of this research was to develop a laboratory test to identify the critical water content "
$d.Content.Find.Execute($merge4, $true, $false, $false, $false, $false, $true, 1, $false, $merge4, 2) | Out-Null

# Paragraph: "A pneumatically-driven device ..." (2 runs -> 1 run, after the spell-checked word)
$merge5 = ". The apparatus applies both compressive and shearing stresses. It may be configured to loading pressures comparable with either youth or professional competition."
$d.Content.Find.Execute($merge5, $true, $false, $false, $false, $false, $true, 1, $false, $merge5, 2) | Out-Null

# Paragraph: "To perform the test ..." (2 runs -> 1 run before the inline math,
# and 2 runs -> 1 run after the inline math)
$merge6a = "To perform the test, a cylindrical soil sample is prepared using Proctor testing equipment and then subjected to wetting and drying cycles. The pneumatic device is actuated to produce several cleat indentations on the soil surface. A 3D scanning technique quantifies the surface" + [char]0x2019 + "s Dirichlet Normal Energy (DNE). "
$d.Content.Find.Execute($merge6a, $true, $false, $false, $false, $false, $true, 1, $false, $merge6a, 2) | Out-Null

$merge6b = " is measured using a combination of 3D scanning and gravimetric methods. The soil is tested at a range of water contents over successive days."
$d.Content.Find.Execute($merge6b, $true, $false, $false, $false, $false, $true, 1, $false, $merge6b, 2) | Out-Null

# Paragraph with the final math + "... to be pinpointed for any soil. Good correspondence ..."
# (3 runs + spellStart/spellEnd proofErr -> 1 run, proofErr removed)
$merge7 = " to be pinpointed for any soil. Good correspondence was achieved across replicate specimens. It is envisaged that the device will find utility in future investigations of infield mix design. "
$d.Content.Find.Execute($merge7, $true, $false, $false, $false, $false, $true, 1, $false, $merge7, 2) | Out-Null

# Switch pStyle "FirstParagraph" -> "BodyText" on paragraphs 4 through 8 (1-based).
for ($i = 4; $i -le 8; $i++) {
    $d.Paragraphs($i).Style = "Body Text"
}

# -----------------------------------------------------------------------
# 3. Style catalogue changes (styles.xml)
# -----------------------------------------------------------------------

# "Body Text" now inherits from "First Paragraph" instead of "Normal".
$bodyTextStyle = $d.Styles("Body Text")
$bodyTextStyle.BaseStyle = "FirstParagraph"

# "First Paragraph" gains explicit Roboto / 10.5pt run formatting.
$firstParaStyle = $d.Styles("First Paragraph")
$firstParaStyle.Font.Name = "Roboto"
$firstParaStyle.Font.Size = 10.5
$firstParaStyle.Font.SizeBi = 10.5

# "Title" gains Roboto for its ascii/hAnsi fonts (east-Asian theme font kept).
$titleStyle = $d.Styles("Title")
$titleStyle.Font.Name = "Roboto"

# "Body Text Char" gains explicit Roboto / Arial(cs) / 10.5pt run formatting.
$bodyTextCharStyle = $d.Styles("Body Text Char")
$bodyTextCharStyle.Font.Name = "Roboto"
$bodyTextCharStyle.Font.Size = 10.5
$bodyTextCharStyle.Font.SizeBi = 10.5
$bodyTextCharStyle.Font.NameBi = "Arial"
